# Updated cryptos list on Thu Nov 28 10:54:31 UTC 2024 with GitHub Actions
#
# Refreshes the Coin/Link/Price/Volume(1h) table on the active sheet with
# freshly scraped values. Two coin pairs also swapped rank position
# (Dai <-> InternetComputer(DFINITY) at rows 31/32, and
# RenderToken <-> Bittensor at rows 38/39), so those rows rewrite all four
# columns (B-E); every other changed row only updates Price (D) and/or
# Volume(1h) (E).
#
# Price cells that look like plain decimals (e.g. "1.45", "0.401") are
# written with a leading apostrophe so Excel stores them as text instead of
# auto-converting to a number -- matching how the sheet already stores every
# other price (including thousand-dot values like "94.951.40") as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText([string]$cellRef, [string]$text) {
    $value = $text
    if ($text.Trim() -match '^-?\d+(\.\d+)?$') {
        # Force text storage for plain-decimal values (prevents Excel from
        # silently turning "1.45" -> 1.45, losing the original text format).
        $value = "'" + $text
    }
    $ws.Range($cellRef).Value = $value
}

Set-CellText "D2" '94.951.40'
Set-CellText "E2" '  +1.31%  '
Set-CellText "D3" '3.607.76'
Set-CellText "E3" '  +4.27%  '
Set-CellText "E4" '  +0.02%  '
Set-CellText "D5" '235.83'
Set-CellText "E5" '  -0.08%  '
Set-CellText "D6" '658.02'
Set-CellText "E6" '  +5.26%  '
Set-CellText "D7" '1.45'
Set-CellText "E7" '  +0.66%  '
Set-CellText "D8" '0.401'
Set-CellText "E8" '  +1.55%  '
Set-CellText "E9" '  -0.02%  '
Set-CellText "D10" '0.990'
Set-CellText "E10" '  -1.69%  '
Set-CellText "D11" '3.607.61'
Set-CellText "E11" '  +4.32%  '
Set-CellText "D12" '42.29'
Set-CellText "E12" '  -3.65%  '
Set-CellText "E13" '  +0.18%  '
Set-CellText "E14" '  +0.38%  '
Set-CellText "D15" '4.297.21'
Set-CellText "E15" '  +4.18%  '
Set-CellText "D16" '94.947.64'
Set-CellText "E16" '  +1.58%  '
Set-CellText "D17" '0.0000252'
Set-CellText "E17" '  +1.12%  '
Set-CellText "D18" '3.608.82'
Set-CellText "E18" '  +3.76%  '
Set-CellText "D19" '7.91'
Set-CellText "E19" '  -5.54%  '
Set-CellText "D20" '12.81'
Set-CellText "E20" '  +9.01%  '
Set-CellText "D21" '17.89'
Set-CellText "E21" '  -2.08%  '
Set-CellText "D22" '3.54'
Set-CellText "E22" '  +4.43%  '
Set-CellText "D23" '0.478'
Set-CellText "E23" '  -7.44%  '
Set-CellText "D24" '503.73'
Set-CellText "E24" '  +0.40%  '
Set-CellText "D25" '0.0000196'
Set-CellText "E25" '  +5.44%  '
Set-CellText "D26" '6.58'
Set-CellText "E26" '  -3.46%  '
Set-CellText "D27" '95.25'
Set-CellText "E27" '  +3.22%  '
Set-CellText "D28" '3.803.92'
Set-CellText "E28" '  +4.10%  '
Set-CellText "D29" '12.48'
Set-CellText "E29" '  +2.08%  '
Set-CellText "D30" '3.12'
Set-CellText "E30" '  +12.95%  '
Set-CellText "B31" 'Dai'
Set-CellText "C31" 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-CellText "D31" '1.00'
Set-CellText "E31" '  -0.05%  '
Set-CellText "B32" 'InternetComputer(DFINITY)'
Set-CellText "C32" 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-CellText "D32" '11.21'
Set-CellText "E32" '  -2.18%  '
Set-CellText "D33" '0.137'
Set-CellText "E33" '  -2.87%  '
Set-CellText "D34" '1.00'
Set-CellText "E34" '  +1.46%  '
Set-CellText "D35" '32.05'
Set-CellText "E35" '  +8.52%  '
Set-CellText "D36" '0.176'
Set-CellText "E36" '  -2.11%  '
Set-CellText "D37" '0.555'
Set-CellText "E37" '  -0.43%  '
Set-CellText "B38" 'RenderToken'
Set-CellText "C38" 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-CellText "D38" '8.06'
Set-CellText "E38" '  +6.33%  '
Set-CellText "B39" 'Bittensor'
Set-CellText "C39" 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
Set-CellText "D39" '565.73'
Set-CellText "E39" '  -0.54%  '
Set-CellText "D40" '1.46'
Set-CellText "E40" '  +2.84%  '
Set-CellText "E41" '  -0.02%  '
Set-CellText "E42" '  -0.23%  '
Set-CellText "D43" '0.910'
Set-CellText "E43" '  -0.67%  '
Set-CellText "D44" '35.38'
Set-CellText "E44" '  +42.65%  '
Set-CellText "E45" '  +0.51%  '
Set-CellText "D46" '23.68'
Set-CellText "D47" '5.58'
Set-CellText "E47" '  +0.62%  '
Set-CellText "D48" '2.22'
Set-CellText "E48" '  +4.34%  '
Set-CellText "D49" '0.0411'
Set-CellText "E49" '  -3.30%  '
Set-CellText "D50" '3.54'
Set-CellText "E50" '  -2.26%  '
Set-CellText "D51" '53.39'
Set-CellText "E51" '  +0.51%  '
